$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: split "Md Nahedul Bar Chowdhury Ruhul 24160119" into three
# runs ("Md ", "Nahedul", " Bar Chowdhury Ruhul 24160119") by toggling
# a formatting property on the middle word and reverting it, which
# forces the run to split at the word boundaries without altering the
# final formatting.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Nahedul") | Out-Null
$rng1.Bold = 1
$rng1.Bold = 0

# ---------------------------------------------------------------------
# Change 2: replace the yellow-highlighted placeholder with the real
# second group member's name/ID, and clear the highlight formatting.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("[Name and ID of other group members]", $false, $false, $false, $false, $false, $true, 1, $false, "Abhishek Ashok Kumar 24140242", 2) | Out-Null
$rng2.HighlightColorIndex = 0

# ---------------------------------------------------------------------
# Change 3: split "Analysis.R code with the appropriate statistics to
# test the hypotheses. " into two runs ("Analysis.R" / " code with the
# appropriate statistics to test the hypotheses. ") the same way.
# ---------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Analysis.R") | Out-Null
$rng3.Bold = 1
$rng3.Bold = 0
